# Commit: "Add files via upload"
#
# The OOXML diff shows a single, tiny change to the speaker notes of the
# 15th slide (ppt/notesSlides/notesSlide5.xml, linked to ppt/slides/slide15.xml):
# the run of text in the "Notes Placeholder 2" shape
#   "=== Results === Chi-square: 245.5 p-value: 0.0000 Degrees of freedom: 27 Cramer's V: 0.153"
# gained a single trailing space
#   "=== Results === Chi-square: 245.5 p-value: 0.0000 Degrees of freedom: 27 Cramer's V: 0.153 "
# (plus PowerPoint's usual dirty="0" bookkeeping attribute that it stamps
# on rPr/endParaRPr after you click into a run and re-save - not something
# content-visible). Re-create that edit through the notes TextRange.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

$notesPage = $s.NotesPage

# Locate the notes body placeholder robustly (by name, falling back to
# the well-known index) rather than assuming shape order.
$notesShape = $null
for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
    $candidate = $notesPage.Shapes.Item($i)
    if ($candidate.Name -eq "Notes Placeholder 2") {
        $notesShape = $candidate
        break
    }
}
if ($notesShape -eq $null) {
    $notesShape = $notesPage.Shapes.Item(2)
}

$notesShape.TextFrame.TextRange.Text = "=== Results === Chi-square: 245.5 p-value: 0.0000 Degrees of freedom: 27 Cramer's V: 0.153 "
